$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at Y:Z, shifting old Y:Z data to AA:AB
$ws.Range("Y:Z").Insert()

# Clear the formatting picked up by the newly inserted cells in row 1 (header row)
$ws.Range("Y1:Z1").Style = "Normal"

# Set new header labels for the two inserted columns
$ws.Range("Y1").Value = "Total Duration"
$ws.Range("Z1").Value = "Current Time"

# Update the row that used to hold the odd "33 mins, 10.0 secs" label/date pair
# (now shifted to AA5/AB5) to the new rhocF parameter row
$ws.Range("AA5").Value = "rhocF"
$ws.Range("AB5").Value = 0.1

# Update selection to match the new active range
$ws.Range("AA2:AB11").Select()
